$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 408, shifting existing rows 408:492 down to 409:493
$ws.Rows(408).Insert()

# Populate the newly inserted row 408 with the new weekly record
$ws.Cells.Item(408, 1).Value = 5
$ws.Cells.Item(408, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(408, 3).Value = "Maule"
$ws.Cells.Item(408, 4).Value = 44995
$ws.Cells.Item(408, 5).Value = 7
$ws.Cells.Item(408, 6).Value = 100114014
$ws.Cells.Item(408, 7).Value = "Betarraga"
$ws.Cells.Item(408, 8).Value = "Sin especificar"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 5000
$ws.Cells.Item(408, 11).Value = 600
$ws.Cells.Item(408, 12).Value = 600
$ws.Cells.Item(408, 13).Value = 600
$ws.Cells.Item(408, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(408, 15).Value = "Región del Maule"
$ws.Cells.Item(408, 16).Value = 120
$ws.Cells.Item(408, 17).Value = 5
$ws.Cells.Item(408, 18).Value = "Hortaliza"

# Preserve the date style (s="2") used by the rest of column D
$ws.Cells.Item(408, 4).NumberFormat = $ws.Cells.Item(409, 4).NumberFormat
